$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values (Price column D, Volume(1h) column E).
# Values are prefixed with a literal apostrophe to force text entry (matches
# the source data which stores these as plain strings, not numbers), then the
# cell style is reset to "Normal" so no stray text-format style is introduced.

$ws.Range("D2").Value = '''63.807.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.69%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.748.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.14%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.20%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''572.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.76%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''156.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.12%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.29%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -1.90%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -4.33%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +0.63%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.380'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.78%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''5.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -17.99%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''3.236.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.03%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''26.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -3.52%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''63.441.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.07%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  -3.45%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.754.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.43%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -0.58%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''4.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -3.54%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''353.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.92%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''6.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -5.32%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.07%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -0.75%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''64.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -3.41%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -1.61%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +0.28%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -3.29%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.0₃0896'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.28%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -5.02%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''6.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.61%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''169.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -3.35%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -5.00%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''20.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -3.28%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +0.25%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -1.80%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -1.51%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -2.56%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -4.11%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''6.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +4.05%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''4.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -4.89%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''327.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -4.60%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''39.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.18%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''21.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.08%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.0582'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -3.47%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''21.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -4.44%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''134.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -3.20%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -3.23%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -4.85%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.51%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +0.16%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''11.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.34%  '
$ws.Range("E51").Style = "Normal"
